# Generate Report for Handoff
# Updates the localization-status report: moves rows that were "low"
# priority / still in "In Translation" to "ht" (handed-off) priority,
# and bumps the relevant handoff timestamps to reflect the newly
# generated handoff.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for rows 4-7
foreach ($r in 4..7) {
    $wsOverview.Range("G$r").Value = "2017-01-03 04:41:27"
}

# zh-cn sheet: Priority (E) + Latest Handoff Datetime (H) for rows 4-7
foreach ($r in 4..7) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2017-01-03 04:41:15"
}

# de-de sheet: Priority (E) for rows 4-7; its "Latest Handoff Datetime"
# (H) shares the same underlying timestamp as the Overview sheet's
# "Latest HO Xliff Generate Date" (G), so bump it to match.
foreach ($r in 4..7) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2017-01-03 04:41:27"
}
